$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10
$ws.Range("D10").Value = "**"

# Row 11
$ws.Range("C11").Value = "***"
$ws.Range("D11").Value = "*"
$ws.Range("G11").Value = "*"

# Row 12
$ws.Range("C12").Value = "***"
$ws.Range("D12").Value = "*"
$ws.Range("G12").Value = "*"

# Row 13
$ws.Range("C13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("H13").Value = "*"

# Row 14
$ws.Range("E14").Value = ""
$ws.Range("G14").Value = "*"

# View changes (zoom + final selection, matching the saved sheetView state)
$ws.Activate() | Out-Null
$window = $excel.ActiveWindow
$window.Zoom = 115
$window.ScrollRow = 13
$window.ScrollColumn = 1
$ws.Range("E26").Select() | Out-Null
